# Weekly update of fruit/vegetable price records (Ramas de apio).
# Applies the cell-level changes described by the commit diff:
#  - updates Fecha (D), Variedad (H), Calidad (I), Volumen (J),
#    Precio minimo/maximo/promedio (K/L/M) and Precio $/Kg (P) on several
#    existing rows
#  - appends a brand-new row 19 with a fresh record

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the existing date-formatted cell (D2) as the template number format so
# that new/edited date cells keep the workbook's original "YYYY-MM-DD HH:MM:SS"
# style instead of Excel inventing a brand-new number format/style entry.
$dateFormat = $ws.Range("D2").NumberFormat

# --- Row 2 ---
$ws.Range("D2").NumberFormat = $dateFormat
$ws.Range("D2").Value = "2022-01-14"
$ws.Range("H2").Value = "Sin especificar"

# --- Row 4 ---
$ws.Range("D4").NumberFormat = $dateFormat
$ws.Range("D4").Value = "2022-05-18"
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 9000
$ws.Range("L4").Value = 9500
$ws.Range("M4").Value = 9250
$ws.Range("P4").Value = 9250

# --- Row 5 ---
$ws.Range("D5").NumberFormat = $dateFormat
$ws.Range("D5").Value = "2021-08-06"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 6000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 6500
$ws.Range("P5").Value = 6500

# --- Row 6 ---
$ws.Range("D6").NumberFormat = $dateFormat
$ws.Range("D6").Value = "2023-01-19"
$ws.Range("J6").Value = 45
$ws.Range("K6").Value = 6000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 6444
$ws.Range("P6").Value = 6444

# --- Row 7 ---
$ws.Range("D7").NumberFormat = $dateFormat
$ws.Range("D7").Value = "2021-03-08"
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7500
$ws.Range("P7").Value = 7500

# --- Row 8 ---
$ws.Range("D8").NumberFormat = $dateFormat
$ws.Range("D8").Value = "2021-02-26"
$ws.Range("H8").Value = "Americana (o)"
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 4000
$ws.Range("L8").Value = 4500
$ws.Range("M8").Value = 4250
$ws.Range("P8").Value = 4250

# --- Row 9 ---
$ws.Range("D9").NumberFormat = $dateFormat
$ws.Range("D9").Value = "2021-10-28"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 5000
$ws.Range("L9").Value = 6000
$ws.Range("M9").Value = 5500
$ws.Range("P9").Value = 5500

# --- Row 10 ---
$ws.Range("D10").NumberFormat = $dateFormat
$ws.Range("D10").Value = "2022-08-16"
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 6000
$ws.Range("M10").Value = 5500
$ws.Range("P10").Value = 5500

# --- Row 11 ---
$ws.Range("D11").NumberFormat = $dateFormat
$ws.Range("D11").Value = "2022-07-22"
$ws.Range("H11").Value = "Americana (o)"

# --- Row 12 ---
$ws.Range("D12").NumberFormat = $dateFormat
$ws.Range("D12").Value = "2022-08-31"
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 5500
$ws.Range("L12").Value = 6000
$ws.Range("M12").Value = 5750
$ws.Range("P12").Value = 5750

# --- Row 14 ---
$ws.Range("D14").NumberFormat = $dateFormat
$ws.Range("D14").Value = "2021-12-29"
$ws.Range("H14").Value = "Americana (o)"
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = 5500
$ws.Range("P14").Value = 5500

# --- Row 15 ---
$ws.Range("D15").NumberFormat = $dateFormat
$ws.Range("D15").Value = "2021-03-04"
$ws.Range("K15").Value = 4000
$ws.Range("L15").Value = 4500
$ws.Range("M15").Value = 4250
$ws.Range("P15").Value = 4250

# --- Row 16 ---
$ws.Range("D16").NumberFormat = $dateFormat
$ws.Range("D16").Value = "2021-12-09"
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 6500
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 6750
$ws.Range("P16").Value = 6750

# --- Row 17 ---
$ws.Range("D17").NumberFormat = $dateFormat
$ws.Range("D17").Value = "2021-03-26"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("K17").Value = 5000
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = 5500
$ws.Range("P17").Value = 5500

# --- Row 18 ---
$ws.Range("D18").NumberFormat = $dateFormat
$ws.Range("D18").Value = "2021-06-24"
$ws.Range("K18").Value = 7000
$ws.Range("L18").Value = 8000
$ws.Range("M18").Value = 7375
$ws.Range("P18").Value = 7375

# --- Row 19 (new row) ---
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").NumberFormat = $dateFormat
$ws.Range("D19").Value = "2022-03-16"
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = 100112017
$ws.Range("G19").Value = "Ramas de apio"
$ws.Range("H19").Value = "Americana (o)"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 60
$ws.Range("K19").Value = 8000
$ws.Range("L19").Value = 9000
$ws.Range("M19").Value = 8500
$ws.Range("N19").Value = "$/atado 7 kilos"
$ws.Range("O19").Value = "Región de Arica y Parinacota"
$ws.Range("P19").Value = 8500
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = "Hortaliza"
